$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (columns B:E) updated
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values (columns B:E) updated
$ws.Range("B2").Value = -4.977969127720586
$ws.Range("C2").Value = 10.722787820974062
$ws.Range("D2").Value = 8.8685582203942523
$ws.Range("E2").Value = 15.404249389654476

# Row 3 data values (columns B:E) updated
$ws.Range("B3").Value = 9.1606185307708188
$ws.Range("C3").Value = 28.120344471525261
$ws.Range("D3").Value = 41.709503591796292
$ws.Range("E3").Value = 9.2166417836742944

# Update the selected range to reflect the narrower area of interest
$ws.Range("B1:E3").Select()
